$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell -> new text value, taken from the updated symbol list.
# Columns D (Price) and E (Volume 1h) hold numeric-looking text that must stay
# text (matching the original inlineStr cells), so we briefly force a Text
# number format while assigning, then restore the "Normal" style so no
# formatting is left behind.
$updates = @(
    @{ Cell = "D2"; Value = '307.28' },
    @{ Cell = "E2"; Value = '-1.34%' },
    @{ Cell = "D3"; Value = '37.38' },
    @{ Cell = "E3"; Value = '-0.81%' },
    @{ Cell = "D4"; Value = '5.128' },
    @{ Cell = "E4"; Value = '0.96%' },
    @{ Cell = "D5"; Value = '0.07804' },
    @{ Cell = "E5"; Value = '0.33%' },
    @{ Cell = "D6"; Value = '4.425' },
    @{ Cell = "E6"; Value = '1.63%' },
    @{ Cell = "D7"; Value = '8.263' },
    @{ Cell = "E7"; Value = '0.45%' },
    @{ Cell = "D8"; Value = '1.876' },
    @{ Cell = "E8"; Value = '-0.41%' },
    @{ Cell = "D9"; Value = '2.997' },
    @{ Cell = "E9"; Value = '5.15%' },
    @{ Cell = "D10"; Value = '0.9277' },
    @{ Cell = "E10"; Value = '0.82%' },
    @{ Cell = "E11"; Value = '-9.23%' },
    @{ Cell = "D12"; Value = '0.1913' },
    @{ Cell = "E12"; Value = '-0.64%' },
    @{ Cell = "D13"; Value = '0.09003' },
    @{ Cell = "E13"; Value = '-3.57%' },
    @{ Cell = "E14"; Value = '-2.69%' },
    @{ Cell = "D15"; Value = '0.09585' },
    @{ Cell = "E15"; Value = '-0.96%' },
    @{ Cell = "D16"; Value = '0.001380' },
    @{ Cell = "E16"; Value = '-0.24%' },
    @{ Cell = "B17"; Value = 'TigerCash' },
    @{ Cell = "C17"; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' },
    @{ Cell = "D17"; Value = '0.005821' },
    @{ Cell = "E17"; Value = '1.18%' },
    @{ Cell = "B18"; Value = 'LEO' },
    @{ Cell = "C18"; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' },
    @{ Cell = "D18"; Value = '3.592' },
    @{ Cell = "E18"; Value = '1.06%' },
    @{ Cell = "B19"; Value = 'BitpandaEcosystemToken' },
    @{ Cell = "C19"; Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best' },
    @{ Cell = "D19"; Value = '0.3474' },
    @{ Cell = "E19"; Value = '2.05%' },
    @{ Cell = "B20"; Value = 'MCDex' },
    @{ Cell = "C20"; Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb' },
    @{ Cell = "D20"; Value = '6.288' },
    @{ Cell = "E20"; Value = '19.38%' },
    @{ Cell = "B21"; Value = 'ProBitToken' },
    @{ Cell = "C21"; Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob' },
    @{ Cell = "D21"; Value = '0.1276' },
    @{ Cell = "E21"; Value = '-0.49%' },
    @{ Cell = "B22"; Value = 'ZBToken' },
    @{ Cell = "C22"; Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb' },
    @{ Cell = "D22"; Value = '0.2587' },
    @{ Cell = "E22"; Value = '-0.09%' },
    @{ Cell = "B23"; Value = 'CoinExToken' },
    @{ Cell = "C23"; Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet' },
    @{ Cell = "D23"; Value = '0.04391' },
    @{ Cell = "E23"; Value = '0.42%' },
    @{ Cell = "D24"; Value = '0.001211' },
    @{ Cell = "E24"; Value = '-0.24%' },
    @{ Cell = "D25"; Value = '0.004250' },
    @{ Cell = "E25"; Value = '-0.23%' },
    @{ Cell = "D26"; Value = '0.0001305' },
    @{ Cell = "E26"; Value = '0.43%' },
    @{ Cell = "D39"; Value = '0.02178' },
    @{ Cell = "E39"; Value = '2.81%' },
    @{ Cell = "D40"; Value = '0.05028' },
    @{ Cell = "E40"; Value = '1.01%' },
    @{ Cell = "D41"; Value = '0.007456' },
    @{ Cell = "E41"; Value = '-2.81%' },
    @{ Cell = "D42"; Value = '0.1347' },
    @{ Cell = "E42"; Value = '0.03%' },
    @{ Cell = "D43"; Value = '0.008723' },
    @{ Cell = "E43"; Value = '-12.04%' },
    @{ Cell = "D44"; Value = '0.002118' },
    @{ Cell = "E44"; Value = '2.82%' },
    @{ Cell = "D45"; Value = '0.007993' },
    @{ Cell = "E45"; Value = '-9.31%' },
    @{ Cell = "D46"; Value = '0.00006569' },
    @{ Cell = "E46"; Value = '-1.36%' },
    @{ Cell = "D47"; Value = '0.00000000752' },
    @{ Cell = "E47"; Value = '0.30%' },
    @{ Cell = "D48"; Value = '0.002866' },
    @{ Cell = "E48"; Value = '-1.58%' },
    @{ Cell = "D49"; Value = '0.001003' },
    @{ Cell = "E49"; Value = '-16.43%' },
    @{ Cell = "D50"; Value = '0.00002106' },
    @{ Cell = "E50"; Value = '0.30%' },
    @{ Cell = "D51"; Value = '0.0002005' },
    @{ Cell = "E51"; Value = '0.30%' }
)

$textColumns = @("D", "E")

foreach ($update in $updates) {
    $cellRef = $update.Cell
    $col = [regex]::Match($cellRef, '^[A-Z]+').Value
    $range = $ws.Range($cellRef)
    if ($textColumns -contains $col) {
        $range.NumberFormat = "@"
        $range.Value = $update.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $update.Value
    }
}
